$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 162, shifting old rows 162-172 down to 164-174.
$ws.Rows.Item(162).Resize(2).Insert()

# New row 162
$ws.Cells.Item(162, 1).Value = 7
$ws.Cells.Item(162, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(162, 3).Value = "Ñuble"
$ws.Cells.Item(162, 4).Value = 44615
$ws.Cells.Item(162, 5).Value = 16
$ws.Cells.Item(162, 6).Value = 100112024
$ws.Cells.Item(162, 7).Value = "Choclo"
$ws.Cells.Item(162, 8).Value = "Choclero"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 20000
$ws.Cells.Item(162, 11).Value = 150
$ws.Cells.Item(162, 12).Value = 180
$ws.Cells.Item(162, 13).Value = 165
$ws.Cells.Item(162, 14).Value = "$/unidad"
$ws.Cells.Item(162, 15).Value = "Región del Maule"
$ws.Cells.Item(162, 16).Value = 165
$ws.Cells.Item(162, 17).Value = 1
$ws.Cells.Item(162, 18).Value = "Hortaliza"

# New row 163
$ws.Cells.Item(163, 1).Value = 7
$ws.Cells.Item(163, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(163, 3).Value = "Ñuble"
$ws.Cells.Item(163, 4).Value = 44615
$ws.Cells.Item(163, 5).Value = 16
$ws.Cells.Item(163, 6).Value = 100112024
$ws.Cells.Item(163, 7).Value = "Choclo"
$ws.Cells.Item(163, 8).Value = "Choclero"
$ws.Cells.Item(163, 9).Value = "Segunda"
$ws.Cells.Item(163, 10).Value = 16000
$ws.Cells.Item(163, 11).Value = 100
$ws.Cells.Item(163, 12).Value = 120
$ws.Cells.Item(163, 13).Value = 110
$ws.Cells.Item(163, 14).Value = "$/unidad"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 110
$ws.Cells.Item(163, 17).Value = 1
$ws.Cells.Item(163, 18).Value = "Hortaliza"

Write-Host "Done"
